$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2025-06-16 -> 2025-06-17, serial 45824 -> 45825) for every data row (2..43).
$ws.Range("C2:C43").Value = 45825
